$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 120.5
$ws.Cells.Item(5, 9).Value = 120.5
$ws.Cells.Item(5, 11).Value = 120.5
$ws.Cells.Item(5, 13).Value = -5.5
$ws.Cells.Item(17, 8).Value = 7866.625
$ws.Cells.Item(17, 10).Value = 7866.625
$ws.Cells.Item(17, 12).Value = 23599.875
$ws.Cells.Item(17, 14).Value = -23935.875
$ws.Cells.Item(40, 8).Value = 5993.4287
$ws.Cells.Item(40, 10).Value = 7030.8
$ws.Cells.Item(40, 12).Value = 7030.8
$ws.Cells.Item(40, 14).Value = -7380.8
$ws.Cells.Item(55, 8).Value = 810.9773
$ws.Cells.Item(55, 9).Value = 330
$ws.Cells.Item(55, 10).Value = 1085.8214
$ws.Cells.Item(55, 11).Value = 330
$ws.Cells.Item(55, 12).Value = 1085.8214
$ws.Cells.Item(55, 13).Value = -116
$ws.Cells.Item(55, 14).Value = -1513.8214
$ws.Cells.Item(88, 8).Value = 303631.7
$ws.Cells.Item(88, 9).Value = 600485
$ws.Cells.Item(88, 11).Value = 600485
$ws.Cells.Item(88, 13).Value = -600079
$ws.Cells.Item(91, 8).Value = 303631.7
$ws.Cells.Item(91, 9).Value = 600485
$ws.Cells.Item(91, 11).Value = 600485
$ws.Cells.Item(91, 13).Value = -599081
$ws.Cells.Item(100, 8).Value = 5110.846
$ws.Cells.Item(100, 9).Value = 3797.6
$ws.Cells.Item(100, 10).Value = 5931.625
$ws.Cells.Item(100, 11).Value = 3797.6
$ws.Cells.Item(100, 12).Value = 5931.625
$ws.Cells.Item(100, 13).Value = -3256.6
$ws.Cells.Item(100, 14).Value = -7013.625
$ws.Cells.Item(112, 8).Value = 2361.1667
$ws.Cells.Item(112, 10).Value = 2361.1667
$ws.Cells.Item(112, 12).Value = 7083.500100000001
$ws.Cells.Item(112, 14).Value = -9299.500100000001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 6385.55
$ws.Cells.Item(2, 9).Value = 1917
$ws.Cells.Item(2, 10).Value = 10041.637
$ws.Cells.Item(2, 11).Value = 1917
$ws.Cells.Item(2, 12).Value = 10041.637
$ws.Cells.Item(2, 13).Value = -1804
$ws.Cells.Item(2, 14).Value = -10267.637
$ws.Cells.Item(4, 8).Value = 249.5
$ws.Cells.Item(4, 9).Value = 249.5
$ws.Cells.Item(4, 11).Value = 249.5
$ws.Cells.Item(4, 13).Value = -133.5
$ws.Cells.Item(5, 8).Value = 253.84616
$ws.Cells.Item(5, 9).Value = 230.1
$ws.Cells.Item(5, 10).Value = 333
$ws.Cells.Item(5, 11).Value = 230.1
$ws.Cells.Item(5, 12).Value = 333
$ws.Cells.Item(5, 13).Value = -118.1
$ws.Cells.Item(5, 14).Value = -557
$ws.Cells.Item(32, 8).Value = 5637.9517
$ws.Cells.Item(32, 9).Value = 5637.9517
$ws.Cells.Item(32, 11).Value = 5637.9517
$ws.Cells.Item(32, 13).Value = -5350.9517
$ws.Cells.Item(61, 8).Value = 7178.952
$ws.Cells.Item(61, 9).Value = 6374.8335
$ws.Cells.Item(61, 10).Value = 12003.667
$ws.Cells.Item(61, 11).Value = 6374.8335
$ws.Cells.Item(61, 12).Value = 12003.667
$ws.Cells.Item(61, 13).Value = -6162.8335
$ws.Cells.Item(61, 14).Value = -12427.667
$ws.Cells.Item(74, 8).Value = 3345
$ws.Cells.Item(74, 9).Value = 2231.6365
$ws.Cells.Item(74, 11).Value = 2231.6365
$ws.Cells.Item(74, 13).Value = -1357.6365
$ws.Cells.Item(77, 8).Value = 3345
$ws.Cells.Item(77, 9).Value = 2231.6365
$ws.Cells.Item(77, 11).Value = 11158.1825
$ws.Cells.Item(77, 13).Value = -6790.182500000001
$ws.Cells.Item(110, 8).Value = 3034
$ws.Cells.Item(110, 9).Value = 2589.1177
$ws.Cells.Item(110, 11).Value = 2589.1177
$ws.Cells.Item(110, 13).Value = -544.1176999999998
$ws.Cells.Item(116, 8).Value = 6385.55
$ws.Cells.Item(116, 9).Value = 1917
$ws.Cells.Item(116, 10).Value = 10041.637
$ws.Cells.Item(116, 11).Value = 1917
$ws.Cells.Item(116, 12).Value = 10041.637
$ws.Cells.Item(116, 13).Value = 377
$ws.Cells.Item(116, 14).Value = -14629.637
$ws.Cells.Item(132, 8).Value = 4196.857
$ws.Cells.Item(132, 10).Value = 12757
$ws.Cells.Item(132, 12).Value = 38271
$ws.Cells.Item(132, 14).Value = -43331
$ws.Cells.Item(136, 8).Value = 7178.952
$ws.Cells.Item(136, 9).Value = 6374.8335
$ws.Cells.Item(136, 10).Value = 12003.667
$ws.Cells.Item(136, 11).Value = 19124.5005
$ws.Cells.Item(136, 12).Value = 36011.001
$ws.Cells.Item(136, 13).Value = -16574.5005
$ws.Cells.Item(136, 14).Value = -41111.001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 6385.55
$ws.Cells.Item(3, 9).Value = 1917
$ws.Cells.Item(3, 10).Value = 10041.637
$ws.Cells.Item(3, 11).Value = 1917
$ws.Cells.Item(3, 12).Value = 10041.637
$ws.Cells.Item(3, 13).Value = -1803
$ws.Cells.Item(3, 14).Value = -10269.637
$ws.Cells.Item(4, 8).Value = 253.84616
$ws.Cells.Item(4, 9).Value = 230.1
$ws.Cells.Item(4, 10).Value = 333
$ws.Cells.Item(4, 11).Value = 230.1
$ws.Cells.Item(4, 12).Value = 333
$ws.Cells.Item(4, 13).Value = -115.1
$ws.Cells.Item(4, 14).Value = -563
$ws.Cells.Item(7, 8).Value = 0
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(20, 8).Value = 7787.6665
$ws.Cells.Item(20, 9).Value = 7545.5
$ws.Cells.Item(20, 10).Value = 8998.5
$ws.Cells.Item(20, 11).Value = 7545.5
$ws.Cells.Item(20, 12).Value = 8998.5
$ws.Cells.Item(20, 13).Value = -7298.5
$ws.Cells.Item(20, 14).Value = -9492.5
$ws.Cells.Item(94, 8).Value = 1386.7455
$ws.Cells.Item(94, 9).Value = 1470.3469
$ws.Cells.Item(94, 11).Value = 1470.3469
$ws.Cells.Item(94, 13).Value = -1019.3469
$ws.Cells.Item(99, 8).Value = 2725.0908
$ws.Cells.Item(99, 9).Value = 1854.1428
$ws.Cells.Item(99, 11).Value = 1854.1428
$ws.Cells.Item(99, 13).Value = -356.1428000000001
$ws.Cells.Item(134, 8).Value = 3235.7917
$ws.Cells.Item(134, 9).Value = 2539.4211
$ws.Cells.Item(134, 11).Value = 7618.263300000001
$ws.Cells.Item(134, 13).Value = -5083.263300000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 5881
$ws.Cells.Item(22, 9).Value = 762.25
$ws.Cells.Item(22, 11).Value = 762.25
$ws.Cells.Item(22, 13).Value = -412.25
$ws.Cells.Item(58, 8).Value = 6085.35
$ws.Cells.Item(58, 9).Value = 3500.1
$ws.Cells.Item(58, 11).Value = 3500.1
$ws.Cells.Item(58, 13).Value = -3297.1
$ws.Cells.Item(136, 8).Value = 6085.35
$ws.Cells.Item(136, 9).Value = 3500.1
$ws.Cells.Item(136, 11).Value = 10500.3
$ws.Cells.Item(136, 13).Value = -7950.299999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(124, 8).Value = 4786.6665
$ws.Cells.Item(124, 10).Value = 5897.1665
$ws.Cells.Item(124, 12).Value = 17691.4995
$ws.Cells.Item(124, 14).Value = -27511.4995

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 7412.407
$ws.Cells.Item(132, 10).Value = 14340.429
$ws.Cells.Item(132, 12).Value = 43021.287
$ws.Cells.Item(132, 14).Value = -48081.287

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 3910.2144
$ws.Cells.Item(22, 9).Value = 1238.4286
$ws.Cells.Item(22, 10).Value = 6582
$ws.Cells.Item(22, 11).Value = 1238.4286
$ws.Cells.Item(22, 12).Value = 6582
$ws.Cells.Item(22, 13).Value = -943.4286
$ws.Cells.Item(22, 14).Value = -7172
$ws.Cells.Item(27, 8).Value = 3910.2144
$ws.Cells.Item(27, 9).Value = 1238.4286
$ws.Cells.Item(27, 10).Value = 6582
$ws.Cells.Item(27, 11).Value = 1238.4286
$ws.Cells.Item(27, 12).Value = 6582
$ws.Cells.Item(27, 13).Value = -1131.4286
$ws.Cells.Item(27, 14).Value = -6796
$ws.Cells.Item(46, 8).Value = 3216.8333
$ws.Cells.Item(46, 9).Value = 1250.5
$ws.Cells.Item(46, 10).Value = 3462.625
$ws.Cells.Item(46, 11).Value = 1250.5
$ws.Cells.Item(46, 12).Value = 3462.625
$ws.Cells.Item(46, 13).Value = -1062.5
$ws.Cells.Item(46, 14).Value = -3838.625
$ws.Cells.Item(55, 8).Value = 1668910
$ws.Cells.Item(55, 9).Value = 3126313
$ws.Cells.Item(55, 10).Value = 3306.6428
$ws.Cells.Item(55, 11).Value = 3126313
$ws.Cells.Item(55, 12).Value = 3306.6428
$ws.Cells.Item(55, 13).Value = -3126140
$ws.Cells.Item(55, 14).Value = -3652.6428
$ws.Cells.Item(61, 8).Value = 13508.833
$ws.Cells.Item(61, 9).Value = 10762
$ws.Cells.Item(61, 10).Value = 19002.5
$ws.Cells.Item(61, 11).Value = 10762
$ws.Cells.Item(61, 12).Value = 19002.5
$ws.Cells.Item(61, 13).Value = -10560
$ws.Cells.Item(61, 14).Value = -19406.5
$ws.Cells.Item(93, 8).Value = 2155.4644
$ws.Cells.Item(93, 9).Value = 2018.875
$ws.Cells.Item(93, 11).Value = 2018.875
$ws.Cells.Item(93, 13).Value = -770.875
$ws.Cells.Item(113, 8).Value = 13508.833
$ws.Cells.Item(113, 9).Value = 10762
$ws.Cells.Item(113, 10).Value = 19002.5
$ws.Cells.Item(113, 11).Value = 10762
$ws.Cells.Item(113, 12).Value = 19002.5
$ws.Cells.Item(113, 13).Value = -8592
$ws.Cells.Item(113, 14).Value = -23342.5
$ws.Cells.Item(122, 8).Value = 99836.09
$ws.Cells.Item(122, 9).Value = 109529.234
$ws.Cells.Item(122, 11).Value = 328587.702
$ws.Cells.Item(122, 13).Value = -326137.702
$ws.Cells.Item(136, 8).Value = 3719.8918
$ws.Cells.Item(136, 9).Value = 2822.5715
$ws.Cells.Item(136, 11).Value = 8467.7145
$ws.Cells.Item(136, 13).Value = -5917.7145

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 5091
$ws.Cells.Item(81, 9).Value = 4012.375
$ws.Cells.Item(81, 10).Value = 7967.3335
$ws.Cells.Item(81, 11).Value = 8024.75
$ws.Cells.Item(81, 12).Value = 15934.667
$ws.Cells.Item(81, 13).Value = -6963.75
$ws.Cells.Item(81, 14).Value = -18056.667
$ws.Cells.Item(84, 8).Value = 5091
$ws.Cells.Item(84, 9).Value = 4012.375
$ws.Cells.Item(84, 10).Value = 7967.3335
$ws.Cells.Item(84, 11).Value = 40123.75
$ws.Cells.Item(84, 12).Value = 79673.33499999999
$ws.Cells.Item(84, 13).Value = -34819.75
$ws.Cells.Item(84, 14).Value = -90281.33499999999
$ws.Cells.Item(107, 8).Value = 1103.4412
$ws.Cells.Item(107, 9).Value = 989.3077
$ws.Cells.Item(107, 11).Value = 2967.9231
$ws.Cells.Item(107, 13).Value = -1047.9231
$ws.Cells.Item(122, 8).Value = 2047.3606
$ws.Cells.Item(122, 9).Value = 1604.2354
$ws.Cells.Item(122, 11).Value = 4812.706200000001
$ws.Cells.Item(122, 13).Value = -2362.706200000001
$ws.Cells.Item(132, 8).Value = 2512.12
$ws.Cells.Item(132, 9).Value = 1833.238
$ws.Cells.Item(132, 10).Value = 6076.25
$ws.Cells.Item(132, 11).Value = 5499.714
$ws.Cells.Item(132, 12).Value = 18228.75
$ws.Cells.Item(132, 13).Value = -2969.714
$ws.Cells.Item(132, 14).Value = -23288.75

# --- Special case: remove BSM!M7 entirely (diff drops the cell) ---
$wsBSM = $wb.Worksheets.Item("BSM")
$wsBSM.Range("M7").ClearContents()

Write-Output "edits applied"